$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Reword the "Kilwins" customer-service bullet point.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Demonstrated customer service while bagging groceries and working the cashier while finishing my degree",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enhanced customer satisfaction through efficient service and problem-solving at the register",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2: Remove the stray empty (formatting-only) paragraph that sits
# directly after the "Demonstrated excellent customer service" bullet,
# right before the ACTIVITIES section heading.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute(
    "Demonstrated excellent customer service",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)

if ($found) {
    $r.Collapse(0)                 # wdCollapseEnd - move to just after the match
    $nextPara = $r.Next(4, 1)      # wdParagraph - range covering the following paragraph mark
    if ($nextPara -ne $null -and $nextPara.Text.Length -le 1) {
        # The paragraph holds no run text (only its end-of-paragraph mark), so
        # this is the empty styling-only paragraph that must be deleted.
        $nextPara.Delete()
    }
}
